# Re-draft the histogram source data in columns A (bin left-edges) and
# B (bin counts). The bin grid grows from 37 rows (A2:B38) to 77 rows
# (A2:B78): every existing row is renumbered onto the new -0.147..0.043
# grid (step 0.0025) and 40 new rows are appended below the old last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The newly appended rows (39:78) need the same formatting as the
# existing bin-edge column (bold/centered/bordered style used by A2:A38).
# Copy the format down before writing values so the style carries over.
$fmtSrc = $ws.Range("A2")
$fmtSrc.Copy()
$ws.Range("A39:A78").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$aVals = @(-0.147,-0.1445,-0.142,-0.1395,-0.137,-0.1345,-0.132,-0.1295,-0.127,-0.1245,-0.122,-0.1195,-0.117,-0.1145,-0.112,-0.1095,-0.107,-0.1045,-0.102,-0.09950000000000001,-0.097,-0.0945,-0.092,-0.0895,-0.08699999999999999,-0.08450000000000001,-0.082,-0.0795,-0.077,-0.0745,-0.07199999999999999,-0.06950000000000001,-0.067,-0.0645,-0.062,-0.0595,-0.057,-0.0545,-0.052,-0.0495,-0.047,-0.0445,-0.042,-0.0395,-0.037,-0.0345,-0.032,-0.0295,-0.027,-0.0245,-0.022,-0.0195,-0.017,-0.0145,-0.012,-0.0095,-0.007,-0.0045,-0.002,0.0005,0.003,0.0055,0.008,0.0105,0.013,0.0155,0.018,0.0205,0.023,0.0255,0.028,0.0305,0.033,0.0355,0.038,0.0405,0.043)
$bVals = @(1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,3,0,0,0,0,0,9,0,0,0,2,0,0,0,3,0,3,0,0,0,1,0,1,0,0,0,1)

$n = $aVals.Length
$data = New-Object 'object[,]' $n,2
for ($i = 0; $i -lt $n; $i++) {
    $data[$i,0] = $aVals[$i]
    $data[$i,1] = $bVals[$i]
}

$ws.Range("A2:B78").Value2 = $data
